$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$cs = $wb.Worksheets.Item("CONVERTION")

# ------------------------------------------------------------------
# 1. Insert a new row at 53 (pushes old rows 53..135 down to 54..136)
# ------------------------------------------------------------------
$ws.Rows("53:53").Insert()

# Bring over the correct cell formatting (borders/number formats/etc.)
# for the freshly inserted row from the row right below it (which is
# the row that used to be row 53 before the insert).
$ws.Range("A54:K54").Copy()
$ws.Range("A53:K53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Grow the table (Table1) so it covers the new last row too.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K136"))

# Restore the calculated-column formulas on the new row (53) and on the
# new final row (136), since they were left stale/blank by the insert.
$ws.Range("G53").Formula  = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G136").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ------------------------------------------------------------------
# 2. Fill in the new "Absence Undertime W/ Pay" entries using the
#    CONVERTION sheet's late calculator (hours in E3, minutes in F3,
#    equivalent value comes back out in G3), exactly like the six new
#    UT(...) rows were produced.
# ------------------------------------------------------------------

# Row 53 -> UT(0-2-35)
$cs.Range("E3").Value = 2
$cs.Range("F3").Value = 35
$ws.Range("B53").Value = "UT(0-2-35)"
$ws.Range("D53").Value = $cs.Range("G3").Value2

# Row 51 -> UT(0-1-48)
$cs.Range("E3").Value = 1
$cs.Range("F3").Value = 48
$ws.Range("B51").Value = "UT(0-1-48)"
$ws.Range("D51").Value = $cs.Range("G3").Value2

# Row 50 -> UT(0-1-50)
$cs.Range("E3").Value = 1
$cs.Range("F3").Value = 50
$ws.Range("B50").Value = "UT(0-1-50)"
$ws.Range("D50").Value = $cs.Range("G3").Value2

# Row 49 -> UT(0-1-34)
$cs.Range("E3").Value = 1
$cs.Range("F3").Value = 34
$ws.Range("B49").Value = "UT(0-1-34)"
$ws.Range("D49").Value = $cs.Range("G3").Value2

# Row 48 -> UT(0-2-29)
$cs.Range("E3").Value = 2
$cs.Range("F3").Value = 29
$ws.Range("B48").Value = "UT(0-2-29)"
$ws.Range("D48").Value = $cs.Range("G3").Value2

# Row 47 -> UT(0-0-39)
$cs.Range("E3").ClearContents()
$cs.Range("F3").Value = 39
$ws.Range("B47").Value = "UT(0-0-39)"
$ws.Range("D47").Value = $cs.Range("G3").Value2

# ------------------------------------------------------------------
# 3. Cosmetic: move the active-cell selection like the author left it.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("F55").Select()
